# Updates the cryptos list worksheet with refreshed price / volume data
# (pulled on Mon Feb 19 15:19:31 UTC 2024) and swaps the TRON / Dogecoin
# rows back to their correct rank order (row 11 <-> row 12).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a price string into a "Price" column cell while keeping it
# stored as text (the source data keeps these numeric-looking strings as
# text, e.g. multi-dot thousands values like "52.179.39"). A leading
# apostrophe forces Excel to treat a plain-looking number as text; we then
# restore the "Normal" style so no stray number-format/quote-prefix style
# is left behind on the cell.
function Set-PriceText($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "52.179.39"
$ws.Range("E2").Value = "  +0.72%  "

# Row 3
$ws.Range("D3").Value = "2.905.90"
$ws.Range("E3").Value = "  +3.63%  "

# Row 5
Set-PriceText $ws.Range("D5") "352.82"
$ws.Range("E5").Value = "  -0.21%  "

# Row 6
Set-PriceText $ws.Range("D6") "114.16"
$ws.Range("E6").Value = "  +1.76%  "

# Row 7
Set-PriceText $ws.Range("D7") "0.556"
$ws.Range("E7").Value = "  -0.25%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
Set-PriceText $ws.Range("D9") "0.621"
$ws.Range("E9").Value = "  -0.85%  "

# Row 10
Set-PriceText $ws.Range("D10") "39.79"
$ws.Range("E10").Value = "  -1.05%  "

# Row 11
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-PriceText $ws.Range("D11") "0.0865"
$ws.Range("E11").Value = "  +2.91%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-PriceText $ws.Range("D12") "0.136"
$ws.Range("E12").Value = "  +0.82%  "

# Row 13
Set-PriceText $ws.Range("D13") "19.81"
$ws.Range("E13").Value = "  -0.62%  "

# Row 14
Set-PriceText $ws.Range("D14") "7.74"
$ws.Range("E14").Value = "  -0.79%  "

# Row 15
$ws.Range("D15").Value = "3.361.78"
$ws.Range("E15").Value = "  +3.54%  "

# Row 16
$ws.Range("D16").Value = "2.897.67"
$ws.Range("E16").Value = "  +3.31%  "

# Row 17
Set-PriceText $ws.Range("D17") "0.985"
$ws.Range("E17").Value = "  +3.99%  "

# Row 18
$ws.Range("D18").Value = "52.231.09"
$ws.Range("E18").Value = "  +0.81%  "

# Row 19
Set-PriceText $ws.Range("D19") "3.34"
$ws.Range("E19").Value = "  +3.31%  "

# Row 20
Set-PriceText $ws.Range("D20") "7.62"
$ws.Range("E20").Value = "  +0.03%  "

# Row 21
Set-PriceText $ws.Range("D21") "14.06"
$ws.Range("E21").Value = "  +3.47%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0977"
$ws.Range("E22").Value = "  +0.10%  "

# Row 23
Set-PriceText $ws.Range("D23") "71.21"
$ws.Range("E23").Value = "  +1.20%  "

# Row 24
Set-PriceText $ws.Range("D24") "269.16"
$ws.Range("E24").Value = "  +0.63%  "

# Row 25
$ws.Range("E25").Value = "  +1.54%  "

# Row 26
$ws.Range("E26").Value = "  +12.50%  "

# Row 27
Set-PriceText $ws.Range("D27") "26.79"
$ws.Range("E27").Value = "  +2.36%  "

# Row 28
Set-PriceText $ws.Range("D28") "1.00"
$ws.Range("E28").Value = "  +0.02%  "

# Row 29
Set-PriceText $ws.Range("D29") "10.65"
$ws.Range("E29").Value = "  +2.49%  "

# Row 30
Set-PriceText $ws.Range("D30") "0.104"
$ws.Range("E30").Value = "  +16.00%  "

# Row 31
Set-PriceText $ws.Range("D31") "6.77"
$ws.Range("E31").Value = "  +10.71%  "

# Row 32
Set-PriceText $ws.Range("D32") "37.59"
$ws.Range("E32").Value = "  -3.82%  "

# Row 33
Set-PriceText $ws.Range("D33") "2.28"
$ws.Range("E33").Value = "  +0.14%  "

# Row 34
Set-PriceText $ws.Range("D34") "6.19"
$ws.Range("E34").Value = "  +11.78%  "

# Row 35
Set-PriceText $ws.Range("D35") "53.14"
$ws.Range("E35").Value = "  +1.92%  "

# Row 36
$ws.Range("E36").Value = "  -1.18%  "

# Row 37
$ws.Range("E37").Value = "  -0.12%  "

# Row 38
Set-PriceText $ws.Range("D38") "3.33"
$ws.Range("E38").Value = "  +4.90%  "

# Row 39
Set-PriceText $ws.Range("D39") "18.85"
$ws.Range("E39").Value = "  -0.63%  "

# Row 40
Set-PriceText $ws.Range("D40") "2.05"
$ws.Range("E40").Value = "  +1.76%  "

# Row 41
Set-PriceText $ws.Range("D41") "2.74"
$ws.Range("E41").Value = "  +8.83%  "

# Row 42
$ws.Range("E42").Value = "  +1.51%  "

# Row 43
Set-PriceText $ws.Range("D43") "23.05"
$ws.Range("E43").Value = "  +4.97%  "

# Row 44
Set-PriceText $ws.Range("D44") "118.77"
$ws.Range("E44").Value = "  -0.70%  "

# Row 45
Set-PriceText $ws.Range("D45") "2.18"
$ws.Range("E45").Value = "  -2.38%  "

# Row 46
Set-PriceText $ws.Range("D46") "2.52"
$ws.Range("E46").Value = "  +1.93%  "

# Row 47
Set-PriceText $ws.Range("D47") "3.52"
$ws.Range("E47").Value = "  -0.08%  "

# Row 48
$ws.Range("D48").Value = "2.179.02"
$ws.Range("E48").Value = "  +3.37%  "

# Row 49
Set-PriceText $ws.Range("D49") "0.262"
$ws.Range("E49").Value = "  +19.57%  "

# Row 50
Set-PriceText $ws.Range("D50") "0.0348"
$ws.Range("E50").Value = "  +11.61%  "

# Row 51
Set-PriceText $ws.Range("D51") "0.954"
$ws.Range("E51").Value = "  -1.93%  "
